$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (U1 / WeMos_D1_mini): mark as DNP
$ws.Range("D3").Value = "DNP"

# Row 8 (J2 / Conn_01x11): mark as DNP
$ws.Range("D8").Value = "DNP"
$ws.Range("A8").Value = "Conn_01x11 - DNP"

$ws.Range("A3").Value = "WeMos_D1_mini - DNP"

$ws.Range("D8").Select()
